$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings are not
# auto-converted to numbers by Excel (the source data is plain text).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.604.48'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '2.548.17'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '301.78'
$ws.Range("E5").Value = '  +1.72%  '
$ws.Range("D6").Value = '97.67'
$ws.Range("D7").Value = '0.573'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("D10").Value = '36.42'
$ws.Range("E10").Value = '  +2.15%  '
$ws.Range("D11").Value = '0.0806'
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("E12").Value = '  +8.67%  '
$ws.Range("E13").Value = '  -1.37%  '
$ws.Range("D14").Value = '2.533.04'
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").Value = '0.877'
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").Value = '14.56'
$ws.Range("E16").Value = '  +3.25%  '
$ws.Range("D17").Value = '42.674.58'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").Value = '13.26'
$ws.Range("E18").Value = '  +6.39%  '
$ws.Range("D19").Value = '0.0₃0982'
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").Value = '6.57'
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("D21").Value = '71.49'
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").Value = '254.51'
$ws.Range("E22").Value = '  -2.12%  '
$ws.Range("E23").Value = '  +2.12%  '
$ws.Range("E24").Value = '  -1.86%  '
$ws.Range("D25").Value = '27.80'
$ws.Range("E25").Value = '  -5.76%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").Value = '10.01'
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("D28").Value = '37.99'
$ws.Range("E28").Value = '  +5.45%  '
$ws.Range("D29").Value = '2.19'
$ws.Range("E29").Value = '  +3.59%  '
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("D31").Value = '155.64'
$ws.Range("E31").Value = '  +3.36%  '
$ws.Range("D32").Value = '2.17'
$ws.Range("E32").Value = '  +0.75%  '
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("D34").Value = '0.0799'
$ws.Range("E34").Value = '  +1.06%  '
$ws.Range("D35").Value = '3.28'
$ws.Range("E35").Value = '  -3.02%  '
$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").Value = '18.38'
$ws.Range("E36").Value = '  +14.25%  '
$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D37").Value = '25.77'
$ws.Range("E37").Value = '  +6.83%  '
$ws.Range("D38").Value = '0.114'
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("E40").Value = '  +33.02%  '
$ws.Range("E41").Value = '  +0.86%  '
$ws.Range("E42").Value = '  -1.84%  '
$ws.Range("D43").Value = '2.072.01'
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("D46").Value = '88.16'
$ws.Range("E46").Value = '  +3.87%  '
$ws.Range("D47").Value = '9.20'
$ws.Range("E47").Value = '  +6.46%  '
$ws.Range("D48").Value = '2.800.82'
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("D49").Value = '74.72'
$ws.Range("E49").Value = '  +8.21%  '
$ws.Range("D50").Value = '102.86'
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("E51").Value = '  +2.09%  '
